# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Update time-to-discovery simulation values (td_sim_1 / record_atd columns)
# for the Radjenovic_2013 IEC dataset, and refresh the average_simulation_TD
# summary row to reflect the corrected td_sim_1 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 3).Value = 30
$ws.Cells.Item(3, 4).Value = 30
$ws.Cells.Item(5, 3).Value = 83
$ws.Cells.Item(5, 4).Value = 83
$ws.Cells.Item(7, 3).Value = 274
$ws.Cells.Item(7, 4).Value = 274
$ws.Cells.Item(9, 3).Value = 43
$ws.Cells.Item(9, 4).Value = 43
$ws.Cells.Item(11, 3).Value = 32
$ws.Cells.Item(11, 4).Value = 32
$ws.Cells.Item(13, 3).Value = 19
$ws.Cells.Item(13, 4).Value = 19
$ws.Cells.Item(15, 3).Value = 172
$ws.Cells.Item(15, 4).Value = 172
$ws.Cells.Item(18, 3).Value = 16
$ws.Cells.Item(18, 4).Value = 16
$ws.Cells.Item(19, 3).Value = 31
$ws.Cells.Item(19, 4).Value = 31
$ws.Cells.Item(21, 3).Value = 264
$ws.Cells.Item(21, 4).Value = 264
$ws.Cells.Item(23, 3).Value = 263
$ws.Cells.Item(23, 4).Value = 263
$ws.Cells.Item(25, 3).Value = 241
$ws.Cells.Item(25, 4).Value = 241
$ws.Cells.Item(27, 3).Value = 95
$ws.Cells.Item(27, 4).Value = 95
$ws.Cells.Item(29, 3).Value = 45
$ws.Cells.Item(29, 4).Value = 45
$ws.Cells.Item(31, 3).Value = 200
$ws.Cells.Item(31, 4).Value = 200
$ws.Cells.Item(33, 3).Value = 20
$ws.Cells.Item(33, 4).Value = 20
$ws.Cells.Item(35, 3).Value = 103
$ws.Cells.Item(35, 4).Value = 103
$ws.Cells.Item(37, 3).Value = 148
$ws.Cells.Item(37, 4).Value = 148
$ws.Cells.Item(38, 3).Value = 1276
$ws.Cells.Item(38, 4).Value = 1276
$ws.Cells.Item(40, 3).Value = 72
$ws.Cells.Item(40, 4).Value = 72
$ws.Cells.Item(42, 3).Value = 386
$ws.Cells.Item(42, 4).Value = 386
$ws.Cells.Item(44, 3).Value = 26
$ws.Cells.Item(44, 4).Value = 26
$ws.Cells.Item(46, 3).Value = 202
$ws.Cells.Item(46, 4).Value = 202
$ws.Cells.Item(48, 3).Value = 112
$ws.Cells.Item(48, 4).Value = 112
$ws.Cells.Item(50, 3).Value = 89
$ws.Cells.Item(50, 4).Value = 89
$ws.Cells.Item(52, 3).Value = 37
$ws.Cells.Item(52, 4).Value = 37
$ws.Cells.Item(54, 3).Value = 187
$ws.Cells.Item(54, 4).Value = 187
$ws.Cells.Item(56, 3).Value = 210
$ws.Cells.Item(56, 4).Value = 210
$ws.Cells.Item(58, 3).Value = 125
$ws.Cells.Item(58, 4).Value = 125
$ws.Cells.Item(60, 3).Value = 230
$ws.Cells.Item(60, 4).Value = 230
$ws.Cells.Item(62, 3).Value = 28
$ws.Cells.Item(62, 4).Value = 28
$ws.Cells.Item(64, 3).Value = 21
$ws.Cells.Item(64, 4).Value = 21
$ws.Cells.Item(66, 3).Value = 203
$ws.Cells.Item(66, 4).Value = 203
$ws.Cells.Item(68, 3).Value = 247
$ws.Cells.Item(68, 4).Value = 247
$ws.Cells.Item(70, 3).Value = 222
$ws.Cells.Item(70, 4).Value = 222
$ws.Cells.Item(72, 3).Value = 209
$ws.Cells.Item(72, 4).Value = 209
$ws.Cells.Item(74, 3).Value = 539
$ws.Cells.Item(74, 4).Value = 539
$ws.Cells.Item(76, 3).Value = 215
$ws.Cells.Item(76, 4).Value = 215
$ws.Cells.Item(78, 3).Value = 106
$ws.Cells.Item(78, 4).Value = 106
$ws.Cells.Item(80, 3).Value = 22
$ws.Cells.Item(80, 4).Value = 22
$ws.Cells.Item(82, 3).Value = 25
$ws.Cells.Item(82, 4).Value = 25
$ws.Cells.Item(84, 3).Value = 218
$ws.Cells.Item(84, 4).Value = 218
$ws.Cells.Item(86, 3).Value = 179
$ws.Cells.Item(86, 4).Value = 179
$ws.Cells.Item(88, 3).Value = 114
$ws.Cells.Item(88, 4).Value = 114
$ws.Cells.Item(90, 3).Value = 126
$ws.Cells.Item(90, 4).Value = 126
$ws.Cells.Item(94, 3).Value = 78
$ws.Cells.Item(94, 4).Value = 78
$ws.Cells.Item(96, 3).Value = 284
$ws.Cells.Item(96, 4).Value = 284

$ws.Cells.Item(97, 3).Value = 174.6875